$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5849801301956177
$ws.Range("B1").Value = 1.220777630805969
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.67427396774292
$ws.Range("E1").Value = 1.443329811096191
